$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name)
$ws.Name = "RoomDisplayData"

# Update cell values:
# A1: "ProfileName" -> "EventsProfileName"
# A2: stays "Test RoomDisplay EventsView" (shared string order changes internally,
#     but the displayed text for A2 is unchanged)
$ws.Range("A1").Value = "EventsProfileName"
$ws.Range("A2").Value = "Test RoomDisplay EventsView"
